$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1021.3913
$ws.Range("J17").Value = 1021.3913
$ws.Range("L17").Value = 3064.1739
$ws.Range("N17").Value = -3400.1739
$ws.Range("H18").Value = 1149.5834
$ws.Range("I18").Value = 708.63635
$ws.Range("K18").Value = 708.63635
$ws.Range("M18").Value = -424.63635
$ws.Range("H33").Value = 659.0571
$ws.Range("I33").Value = 158.17241
$ws.Range("J33").Value = 3080
$ws.Range("K33").Value = 158.17241
$ws.Range("L33").Value = 3080
$ws.Range("M33").Value = 70.82758999999999
$ws.Range("N33").Value = -3538
$ws.Range("H116").Value = 3808
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 5013.3335
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 5013.3335
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -11897.3335
$ws.Range("H132").Value = 6585137.5
$ws.Range("I132").Value = 6950701
$ws.Range("K132").Value = 20852103
$ws.Range("M132").Value = -20849573

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5943.013
$ws.Range("I32").Value = 4683.757
$ws.Range("J32").Value = 37004.668
$ws.Range("K32").Value = 4683.757
$ws.Range("L32").Value = 37004.668
$ws.Range("M32").Value = -4396.757
$ws.Range("N32").Value = -37578.668
$ws.Range("H74").Value = 865.9400000000001
$ws.Range("I74").Value = 818.6279
$ws.Range("J74").Value = 1156.5714
$ws.Range("K74").Value = 818.6279
$ws.Range("L74").Value = 1156.5714
$ws.Range("M74").Value = 55.37210000000005
$ws.Range("N74").Value = -2904.5714
$ws.Range("H77").Value = 865.9400000000001
$ws.Range("I77").Value = 818.6279
$ws.Range("J77").Value = 1156.5714
$ws.Range("K77").Value = 4093.1395
$ws.Range("L77").Value = 5782.857
$ws.Range("M77").Value = 274.8605000000002
$ws.Range("N77").Value = -14518.857
$ws.Range("H110").Value = 91010360
$ws.Range("I110").Value = 143015280
$ws.Range("J110").Value = 1741.5
$ws.Range("K110").Value = 143015280
$ws.Range("L110").Value = 1741.5
$ws.Range("M110").Value = -143013235
$ws.Range("N110").Value = -5831.5
$ws.Range("H132").Value = 2047.4791
$ws.Range("J132").Value = 2998.3333
$ws.Range("L132").Value = 8994.999899999999
$ws.Range("N132").Value = -14054.9999
$ws.Range("H135").Value = 44192.57
$ws.Range("J135").Value = 44192.57
$ws.Range("L135").Value = 44192.57
$ws.Range("N135").Value = -54332.57

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 64630.562
$ws.Range("I20").Value = 85599.086
$ws.Range("J20").Value = 1725
$ws.Range("K20").Value = 85599.086
$ws.Range("L20").Value = 1725
$ws.Range("M20").Value = -85352.086
$ws.Range("N20").Value = -2219
$ws.Range("H22").Value = 335
$ws.Range("I22").Value = 335
$ws.Range("K22").Value = 335
$ws.Range("M22").Value = -162
$ws.Range("H33").Value = 12600
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 12600
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 12600
$ws.Range("N33").Value = -13272
$ws.Range("H86").Value = 46473.68
$ws.Range("I86").Value = 67499.17999999999
$ws.Range("J86").Value = 1794.5
$ws.Range("K86").Value = 67499.17999999999
$ws.Range("L86").Value = 1794.5
$ws.Range("M86").Value = -66376.17999999999
$ws.Range("N86").Value = -4040.5
$ws.Range("H89").Value = 46473.68
$ws.Range("I89").Value = 67499.17999999999
$ws.Range("J89").Value = 1794.5
$ws.Range("K89").Value = 337495.9
$ws.Range("L89").Value = 8972.5
$ws.Range("M89").Value = -331879.9
$ws.Range("N89").Value = -20204.5
$ws.Range("H107").Value = 142858500
$ws.Range("I107").Value = 333333700
$ws.Range("J107").Value = 2086
$ws.Range("K107").Value = 333333700
$ws.Range("L107").Value = 2086
$ws.Range("M107").Value = -333331780
$ws.Range("N107").Value = -5926
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N109").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 377.125
$ws.Range("I22").Value = 205.4
$ws.Range("J22").Value = 663.3333
$ws.Range("K22").Value = 205.4
$ws.Range("L22").Value = 663.3333
$ws.Range("M22").Value = 144.6
$ws.Range("N22").Value = -1363.3333
$ws.Range("H86").Value = 3521.2856
$ws.Range("I86").Value = 3275
$ws.Range("J86").Value = 3849.6667
$ws.Range("K86").Value = 3275
$ws.Range("L86").Value = 3849.6667
$ws.Range("M86").Value = -2152
$ws.Range("N86").Value = -6095.6667
$ws.Range("H89").Value = 3521.2856
$ws.Range("I89").Value = 3275
$ws.Range("J89").Value = 3849.6667
$ws.Range("K89").Value = 16375
$ws.Range("L89").Value = 19248.3335
$ws.Range("M89").Value = -10759
$ws.Range("N89").Value = -30480.3335
$ws.Range("H99").Value = 12136
$ws.Range("I99").Value = 3786
$ws.Range("J99").Value = 26748.5
$ws.Range("K99").Value = 3786
$ws.Range("L99").Value = 26748.5
$ws.Range("M99").Value = -2288
$ws.Range("N99").Value = -29744.5
$ws.Range("H126").Value = 12136
$ws.Range("I126").Value = 3786
$ws.Range("J126").Value = 26748.5
$ws.Range("K126").Value = 11358
$ws.Range("L126").Value = 80245.5
$ws.Range("M126").Value = -8888
$ws.Range("N126").Value = -85185.5
$ws.Range("H134").Value = 1702.8334
$ws.Range("I134").Value = 1581.5555
$ws.Range("J134").Value = 2066.6667
$ws.Range("K134").Value = 4744.666499999999
$ws.Range("L134").Value = 6200.000100000001
$ws.Range("M134").Value = -2209.666499999999
$ws.Range("N134").Value = -11270.0001
$ws.Range("H140").Value = 54750
$ws.Range("J140").Value = 54750
$ws.Range("L140").Value = 54750
$ws.Range("N140").Value = -65110

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1086.862
$ws.Range("I5").Value = 924.44446
$ws.Range("J5").Value = 1352.6364
$ws.Range("K5").Value = 2773.33338
$ws.Range("L5").Value = 4057.9092
$ws.Range("M5").Value = -2661.33338
$ws.Range("N5").Value = -4281.9092
$ws.Range("H121").Value = 1053010.1
$ws.Range("I121").Value = 3117.75
$ws.Range("J121").Value = 1227992.1
$ws.Range("K121").Value = 9353.25
$ws.Range("L121").Value = 3683976.3
$ws.Range("M121").Value = -8043.25
$ws.Range("N121").Value = -3686596.3
$ws.Range("H122").Value = 663.25
$ws.Range("J122").Value = 699.6667
$ws.Range("L122").Value = 6297.0003
$ws.Range("N122").Value = -11197.0003
$ws.Range("H131").Value = 788.12
$ws.Range("J131").Value = 804.24744
$ws.Range("L131").Value = 2412.74232
$ws.Range("N131").Value = -12492.74232
$ws.Range("H133").Value = 3000
$ws.Range("I133").Value = 3000
$ws.Range("K133").Value = 9000
$ws.Range("M133").Value = -3940
$ws.Range("H135").Value = 1086.862
$ws.Range("I135").Value = 924.44446
$ws.Range("J135").Value = 1352.6364
$ws.Range("K135").Value = 8320.00014
$ws.Range("L135").Value = 12173.7276
$ws.Range("M135").Value = -5785.00014
$ws.Range("N135").Value = -17243.7276
$ws.Range("H138").Value = 2407.375
$ws.Range("I138").Value = 1469.75
$ws.Range("J138").Value = 3345
$ws.Range("K138").Value = 4409.25
$ws.Range("L138").Value = 10035
$ws.Range("M138").Value = 730.75
$ws.Range("N138").Value = -20315

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1629.9
$ws.Range("I113").Value = 1300
$ws.Range("J113").Value = 1666.5555
$ws.Range("K113").Value = 1300
$ws.Range("L113").Value = 1666.5555
$ws.Range("M113").Value = 870
$ws.Range("N113").Value = -6006.5555
$ws.Range("H126").Value = 5885288.5
$ws.Range("I126").Value = 5397.3335
$ws.Range("J126").Value = 8405242
$ws.Range("K126").Value = 16192.0005
$ws.Range("L126").Value = 25215726
$ws.Range("M126").Value = -13722.0005
$ws.Range("N126").Value = -25220666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4461.75
$ws.Range("I7").Value = 2668
$ws.Range("K7").Value = 2668
$ws.Range("M7").Value = -2556
$ws.Range("H126").Value = 4461.75
$ws.Range("I126").Value = 2668
$ws.Range("K126").Value = 8004
$ws.Range("M126").Value = -5534
$ws.Range("H129").Value = 34920
$ws.Range("J129").Value = 34920
$ws.Range("L129").Value = 34920
$ws.Range("N129").Value = -44920

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 200755
$ws.Range("J107").Value = 200755
$ws.Range("L107").Value = 602265
$ws.Range("N107").Value = -606105
$ws.Range("H126").Value = 1464.2858
$ws.Range("I126").Value = 1277.7778
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 3833.3334
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -1363.3334
$ws.Range("N126").Value = -10340
$ws.Range("H129").Value = 39790
$ws.Range("J129").Value = 39790
$ws.Range("L129").Value = 39790
$ws.Range("N129").Value = -49790
